# Update "想去人数" (column F) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Row = 2;  Value = 4 },
    @{ Row = 4;  Value = 75 },
    @{ Row = 7;  Value = 1648 },
    @{ Row = 11; Value = 1526 },
    @{ Row = 15; Value = 253 },
    @{ Row = 16; Value = 191 },
    @{ Row = 21; Value = 270 },
    @{ Row = 22; Value = 152 },
    @{ Row = 24; Value = 210 }
)

$ws1 = $wb.Worksheets.Item("展览")
foreach ($u in $updates) {
    $ws1.Cells.Item($u.Row, 6).Value = $u.Value
}

$updates2 = @(
    @{ Row = 2;  Value = 4 },
    @{ Row = 4;  Value = 75 },
    @{ Row = 7;  Value = 1648 },
    @{ Row = 12; Value = 1526 },
    @{ Row = 16; Value = 253 },
    @{ Row = 17; Value = 191 },
    @{ Row = 22; Value = 270 },
    @{ Row = 23; Value = 152 },
    @{ Row = 25; Value = 210 }
)

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($u in $updates2) {
    $ws4.Cells.Item($u.Row, 6).Value = $u.Value
}
